$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(1.199071582251861,1.048591397316329,0.9558088060554724,0.917904714957956,0.9116051525190301,0.9552979966888415,1.14726806545184,1.520539224725894,1.79272185520415,1.916074158366541,1.962715104956658,1.952673299130709,1.919912755160908,1.899836784276829,1.784650876091746,1.713866800003302,1.673110088650446,1.659303174080435,1.721406418501203,1.929537241980029,2.065153987059546,1.992811176225359,1.717997949881578,1.419912634236312)
$colC = @(0.3178250702425487,0.2796713744100998,0.256141874402374,0.2465283047361595,0.244930485682147,0.2560123230507543,0.3046915072031879,0.3993042774215496,0.4682666605880286,0.4995135566634303,0.5113273650162,0.5087838952360357,0.5004858654886561,0.4954006215055529,0.4662220282388034,0.4482894495336041,0.4379634403277919,0.4344652527788639,0.4501996163149329,0.5029237123263783,0.5372725085346701,0.5189501926435582,0.4493360807224462,0.3738031615654336)
$colE = @(0.8798451468993562,0.8482515784806566,0.8292103933822119,0.8215409155136939,0.8202728406842397,0.8291065957569685,0.8688775488406151,0.9497050979195905,1.010827414165973,1.039013438203341,1.049741612895303,1.047428675106843,1.039894954942952,1.035287460357381,1.008993055048222,0.9929598635472843,0.9837738780308314,0.9806698262874392,0.9946629110305736,1.042106306133874,1.073432338126736,1.056683877019253,0.9938928644848346,0.9275344979854339)
$colF = @(1.88625115010565,1.858967255064229,1.843517162896703,1.837547803209446,1.836576302764101,1.84343533638409,1.876572866663764,1.951931273305348,2.013690099286435,2.043189132929598,2.05456268570768,2.052104148421449,2.04412077005469,2.039257172530171,2.011790611802184,1.995301183207957,1.98594905527132,1.982805255151604,1.997042826771221,2.046460167522156,2.079940397169963,2.061962804515247,1.996255031491316,1.930427050121679)
$colG = @(0.1612604381511389,0.1650523209733024,0.1676665267281336,0.1688033102836286,0.1689963763168656,0.1676815689757802,0.1625082632025467,0.1546520295998803,0.1503019923563684,0.1486375191403013,0.1480528998865012,0.1481767689270939,0.1485885039499593,0.1488466665217558,0.1504171362110043,0.1514614259868452,0.1520916366429645,0.1523100781425555,0.1513471964871087,0.1484663238715171,0.1468499872374025,0.1476881176386229,0.1513987467470059,0.1565295687161097)
$colH = @(0.3483257670118078,0.3543665596727053,0.3583439427460675,0.3600321044057111,0.360316486260416,0.3583664373684776,0.350352870568301,0.3367725966994755,0.3281034506119127,0.3244452388516734,0.3231011467799405,0.323388786594677,0.324333833156949,0.324918070482326,0.3283482796327704,0.3305258126905244,0.331805128858953,0.3322428914528288,0.3302912299229419,0.3240551309089668,0.3202196551762597,0.3222446953678499,0.3303971993797887,0.3402171889872747)
$colI = @(0.2169278018084109,0.2240789126147897,0.2287472014286971,0.2307190987788095,0.2310507232801176,0.2287735139159999,0.21933576953136,0.2030405858959998,0.1924311623141364,0.1879035396800344,0.1862322524509592,0.1865902671678259,0.1877651739798423,0.1884904757965038,0.1927330866654819,0.195412482529985,0.1969816945379712,0.1975178195109701,0.1951243465853412,0.1874189000279713,0.182635024240132,0.1851651114298569,0.1952545231822791,0.207210584299911)
$colJ = @(0.01794134904406164,0.01834828464058802,0.0186208871742064,0.01873766846595526,0.0187574030470774,0.01862243910806072,0.01807692512259607,0.01718874681567861,0.01664842401051381,0.01642731281276788,0.01634716055871444,0.01636426321017659,0.01642064677390564,0.01645565010023731,0.01666337326892986,0.01679714701274726,0.01687641185750266,0.01690364747184603,0.01678266607558321,0.01640398823044364,0.01617736734554676,0.01629640053966774,0.01678920556929775,0.0174094257581956)
$colO = @(0.9206544587897696,0.941088422014964,0.9547946878857374,0.960670535704935,0.9616637222684759,0.9548727572538809,0.9274586743551225,0.8829545543946153,0.8559713652977337,0.8449521275111209,0.840961248238969,0.8418126459281865,0.8446201430369484,0.8463635396496585,0.8567168762285178,0.863390829980375,0.8673476112761875,0.8687075528765433,0.8626681437645374,0.8437905667597079,0.832513648818761,0.8384348654437872,0.862994497016345,0.8939958520573157)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Range("B" + $row).Value = $colB[$i]
    $ws.Range("C" + $row).Value = $colC[$i]
    $ws.Range("E" + $row).Value = $colE[$i]
    $ws.Range("F" + $row).Value = $colF[$i]
    $ws.Range("G" + $row).Value = $colG[$i]
    $ws.Range("H" + $row).Value = $colH[$i]
    $ws.Range("I" + $row).Value = $colI[$i]
    $ws.Range("J" + $row).Value = $colJ[$i]
    $ws.Range("O" + $row).Value = $colO[$i]
}

Write-Output "Applied 380 kV case update"